$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row entry optionally carries B (Coin name), C (Link), D (Price, stored
# as text in the sheet) and E (Volume 1h %, stored as text). Only the keys
# present are written.
$rows = @(
    @{ R=2;  D="98.285.57";  E="  -0.41%  " }
    @{ R=3;  D="3.425.49";   E="  +2.68%  " }
    @{ R=4;  E="  -0.02%  " }
    @{ R=5;  D="258.04";     E="  +0.68%  " }
    @{ R=6;  D="658.54";     E="  +2.57%  " }
    @{ R=7;  D="1.48";       E="  -5.78%  " }
    @{ R=8;  D="0.443";      E="  +3.49%  " }
    @{ R=9;  D="1.07";       E="  -1.60%  " }
    @{ R=10; E="  +0.00%  " }
    @{ R=11; D="3.423.12";   E="  +2.73%  " }
    @{ R=12; E="  +4.10%  " }
    @{ R=13; D="42.47";      E="  -2.47%  " }
    @{ R=14; E="  +15.57%  " }
    @{ R=15; D="0.0000271"; E="  -0.16%  " }
    @{ R=16; D="98.054.15";  E="  -0.39%  " }
    @{ R=17; D="4.063.61";   E="  +2.38%  " }
    @{ R=18; D="9.41";       E="  +32.86%  " }
    @{ R=19; D="0.597";      E="  +35.16%  " }
    @{ R=20; D="3.427.78";   E="  +2.73%  " }
    @{ R=21; D="17.99";      E="  +8.42%  " }
    @{ R=22; B="Uniswap";    C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni";         D="10.93"; E="  +8.28%  " }
    @{ R=23; B="SuiNetwork"; C="https://coinranking.com/coin/3xJluUMvp+suinetwork-sui";       D="3.51";  E="  -0.40%  " }
    @{ R=24; D="518.12";     E="  -3.90%  " }
    @{ R=25; D="0.0000208"; E="  +2.54%  " }
    @{ R=26; D="6.42";       E="  +4.27%  " }
    @{ R=27; D="101.94";     E="  +1.05%  " }
    @{ R=28; D="13.19";      E="  +5.63%  " }
    @{ R=29; D="3.610.50";   E="  +2.74%  " }
    @{ R=30; D="0.155";      E="  +3.10%  " }
    @{ R=31; D="12.05";      E="  +9.91%  " }
    @{ R=32; D="0.200";      E="  +5.62%  " }
    @{ R=33; D="0.998";      E="  -0.20%  " }
    @{ R=34; D="0.589";      E="  +12.73%  " }
    @{ R=35; E="  +0.22%  " }
    @{ R=36; E="  +15.14%  " }
    @{ R=37; D="30.21";      E="  +3.68%  " }
    @{ R=38; D="7.96";       E="  +4.86%  " }
    @{ R=39; B="Bittensor";  C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao";        D="541.18"; E="  +3.87%  " }
    @{ R=40; B="Fetch.AI";   C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet";      D="1.47";  E="  +12.09%  " }
    @{ R=41; E="  +0.92%  " }
    @{ R=42; E="  +0.06%  " }
    @{ R=43; D="9.46";       E="  +22.36%  " }
    @{ R=44; B="ARBITRUM";   C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb";         D="0.892"; E="  +10.05%  " }
    @{ R=45; B="Filecoin";   C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil";         D="6.09";  E="  +24.29%  " }
    @{ R=46; D="24.74";      E="  +0.04%  " }
    @{ R=47; B="VeChain";    C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet";      D="0.0440"; E="  +12.36%  " }
    @{ R=48; D="3.71";       E="  -3.72%  " }
    @{ R=49; E="  +4.38%  " }
    @{ R=50; D="1.68";       E="  +14.56%  " }
    @{ R=51; D="2.12";       E="  +4.51%  " }
)

foreach ($row in $rows) {
    $r = $row.R

    if ($row.ContainsKey("B")) {
        $ws.Cells.Item($r, 2).Value = $row.B
    }
    if ($row.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $row.C
    }
    if ($row.ContainsKey("D")) {
        # The Price column stores plain-looking numbers as TEXT in the
        # original workbook (t="inlineStr"). Assigning the bare string
        # directly would let the numeric-looking ones ("258.04", "0.0440",
        # ...) get auto-coerced into real numbers, losing formatting
        # (trailing zeros) and the original text cell type. Force text by
        # flipping the cell to the "@" (Text) number format only for the
        # instant of the write, then restore the default "Normal" style so
        # the cell's style index is unchanged from the original (s absent).
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $row.D
        $cell.Style = "Normal"
    }
    if ($row.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $row.E
    }
}
